$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the rows for the 4 discontinued models.
#    Delete bottom-to-top so earlier row numbers stay valid while we work.
#    Row 20 -> T109.410.11.072.00
#    Row 19 -> T41.1.123.57
#    Row 11 -> T126.010.11.013.00
#    Row 9  -> T120.417.11.041.01
$ws.Rows(20).Delete()
$ws.Rows(19).Delete()
$ws.Rows(11).Delete()
$ws.Rows(9).Delete()

# 2. Update "去處" (column C) for every remaining data row: the vendor code
#    was renamed from C106 (SAHS) to C107 (SAHS).
$ws.Range("C2:C16").Value = "C107 (SAHS) 0.60"

# 3. Update "數量" (column B, quantity) for the remaining rows to their new counts.
$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 10
$ws.Range("B4").Value = 10
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 5
$ws.Range("B9").Value = 5
$ws.Range("B10").Value = 5
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 5
$ws.Range("B13").Value = 1
$ws.Range("B14").Value = 3
$ws.Range("B15").Value = 3
$ws.Range("B16").Value = 5

# 4. Restore the sheet view: no frozen/scrolled top-left cell, selection on C20.
$ws.Range("C20").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
